# Applies the cryptos list refresh: updates Coin/Link/Price/Volume for rows 2-51
# on the active sheet. Row 9 onward shift up one coin (OKB drops off,
# Decentraland is appended at row 51); D/E hold fresh price & volume text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ R = 2; B = "Bitcoin"; C = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D = "25.846.20"; E = "  -0.26%  " },
    @{ R = 3; B = "Ethereum"; C = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D = "1.737.15"; E = "  -0.33%  " },
    @{ R = 4; B = "TetherUSD"; C = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D = "1.000"; E = "  +0.08%  " },
    @{ R = 5; B = "BNB"; C = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D = "240.80"; E = "  +3.91%  " },
    @{ R = 6; B = "USDC"; C = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D = "1.001"; E = "  +0.09%  " },
    @{ R = 7; B = "XRP"; C = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D = "0.5183"; E = "  -1.22%  " },
    @{ R = 8; B = "Cardano"; C = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D = "0.2741"; E = "  -0.95%  " },
    @{ R = 9; B = "Dogecoin"; C = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D = "0.06150"; E = "  +0.18%  " },
    @{ R = 10; B = "WrappedEther"; C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D = "1.740.21"; E = "  -0.12%  " },
    @{ R = 11; B = "TRON"; C = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D = "0.07167"; E = "  +0.71%  " },
    @{ R = 12; B = "Solana"; C = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D = "14.92"; E = "  -2.37%  " },
    @{ R = 13; B = "Polygon"; C = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D = "0.6409"; E = "  -0.76%  " },
    @{ R = 14; B = "Polkadot"; C = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D = "4.607"; E = "  +1.69%  " },
    @{ R = 15; B = "Litecoin"; C = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D = "77.20"; E = "  -0.13%  " },
    @{ R = 16; B = "Dai"; C = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; D = "1.000"; E = "  +0.10%  " },
    @{ R = 17; B = "BinanceUSD"; C = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D = "1.001"; E = "  +0.11%  " },
    @{ R = 18; B = "WrappedBTC"; C = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D = "25.874.41"; E = "  -0.05%  " },
    @{ R = 19; B = "Avalanche"; C = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D = "11.70"; E = "  +1.11%  " },
    @{ R = 20; B = "ShibaInu"; C = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D = "0.000006768"; E = "  +1.13%  " },
    @{ R = 21; B = "WrappedliquidstakedEther2.0"; C = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D = "1.964.00"; E = "  +0.23%  " },
    @{ R = 22; B = "Uniswap"; C = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D = "4.258"; E = "  -0.36%  " },
    @{ R = 23; B = "Cosmos"; C = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D = "8.626"; E = "  -1.85%  " },
    @{ R = 24; B = "Chainlink"; C = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D = "5.253"; E = "  +1.34%  " },
    @{ R = 25; B = "Monero"; C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D = "138.06"; E = "  -1.55%  " },
    @{ R = 26; B = "Toncoin"; C = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D = "1.520"; E = "  -0.03%  " },
    @{ R = 27; B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "15.19"; E = "  -0.16%  " },
    @{ R = 28; B = "LidoDAOToken"; C = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D = "1.766"; E = "  -2.27%  " },
    @{ R = 29; B = "BitcoinCash"; C = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D = "104.92"; E = "  +2.29%  " },
    @{ R = 30; B = "InternetComputer(DFINITY)"; C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D = "3.942"; E = "  +5.25%  " },
    @{ R = 31; B = "Stellar"; C = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D = "0.08256"; E = "  -1.19%  " },
    @{ R = 32; B = "Filecoin"; C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D = "3.658"; E = "  +2.45%  " },
    @{ R = 33; B = "Hedera"; C = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; D = "0.04629"; E = "  +2.19%  " },
    @{ R = 34; B = "HuobiToken"; C = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D = "2.644"; E = "  +1.33%  " },
    @{ R = 35; B = "ARBITRUM"; C = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D = "0.9861"; E = "  +0.64%  " },
    @{ R = 36; B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "0.6154"; E = "  -1.51%  " },
    @{ R = 37; B = "MXToken"; C = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D = "2.692"; E = "  -0.17%  " },
    @{ R = 38; B = "VeChain"; C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D = "0.01598"; E = "  +0.41%  " },
    @{ R = 39; B = "RenderToken"; C = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D = "1.917"; E = "  -0.70%  " },
    @{ R = 40; B = "PaxDollar"; C = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"; D = "1.000"; E = "  +0.10%  " },
    @{ R = 41; B = "Quant"; C = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D = "99.56"; E = "  -0.82%  " },
    @{ R = 42; B = "TheSandbox"; C = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D = "0.3837"; E = "  -1.05%  " },
    @{ R = 43; B = "TrustWalletToken"; C = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D = "0.7446"; E = "  +0.49%  " },
    @{ R = 44; B = "FraxShare"; C = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D = "4.991"; E = "  -0.94%  " },
    @{ R = 45; B = "Algorand"; C = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D = "0.1123"; E = "  -0.39%  " },
    @{ R = 46; B = "Aptos"; C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D = "6.234"; E = "  -0.31%  " },
    @{ R = 47; B = "Cronos"; C = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D = "0.05242"; E = "  -1.81%  " },
    @{ R = 48; B = "Aave"; C = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D = "54.73"; E = "  +1.80%  " },
    @{ R = 49; B = "Elrond"; C = "https://coinranking.com/coin/omwkOTglq+elrond-egld"; D = "30.49"; E = "  +0.91%  " },
    @{ R = 50; B = "EnergySwap"; C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D = "7.596"; E = "  -0.75%  " },
    @{ R = 51; B = "Decentraland"; C = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; D = "0.3407"; E = "  -0.91%  " }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    # Price column holds text that often LOOKS numeric ("1.000", "240.80", ...).
    # Force text entry, then clear the Text number-format stamp so the cell
    # keeps its original (un-styled) look, matching the source data.
    $ws.Range("D$r").NumberFormat = "@"
    $ws.Range("D$r").Value = $row.D
    $ws.Range("D$r").ClearFormats()
    $ws.Range("E$r").Value = $row.E
}